# Generate Report for Handback
# Adds a new handback entry (deacab6e-5088-41e0-b5b3-4ecc97525912) as row 4
# to the "Overview", "zh-cn" and "de-de" worksheets, mirroring the existing
# two rows already present on each sheet.

$wb = $excel.ActiveWorkbook

$uuid      = "deacab6e-5088-41e0-b5b3-4ecc97525912"
$fileHash  = "d539915c154dd968b6735e46ff6a28bcd3d844ff"
$status    = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d50dfab98c3a521bf0c0407ecd061bd64381aace/e2e/$uuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $overviewUrl, $null, $null, "$uuid.md")
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status

# ---------------------------------------------------------------------
# Sheets "zh-cn" and "de-de": same 12-column layout
#   A Source File Name | B File Extension | C Status
#   D Correspond Handoff File | E Correspond Handoff Datetime
#   F Target File | G Correspond Handback File | H Correspond Handback DateTime
#   I Reference Tokens | J Handoff Reason | K Dependency From | L Error Detail
# ---------------------------------------------------------------------
$langs = @(
    @{ Name = "zh-cn"; Fly = "zhcn"; HandoffTime = "2016-03-31 06:53:51"; HandbackTime = "2016-03-31 06:54:48";
       HandoffSha = "9582be2008c481443713618aeea5651a539c25c2"; TargetSha = "a20e349ccf0ee128a38d086dc42f849410750f33"; HandbackSha = "9d92af9190885188f9ef16834818cf5a33c5f27c" },
    @{ Name = "de-de"; Fly = "dede"; HandoffTime = "2016-03-31 06:54:03"; HandbackTime = "2016-03-31 06:55:06";
       HandoffSha = "0f6e4b0954a27b6c4b603b71c68bac15917ce022"; TargetSha = "abf9e811a934b45b8ad7818d6666d45dcdfcb288"; HandbackSha = "0bd250e665c1c2c32f8fbcf4c0ad07a6db2a7734" }
)

foreach ($lang in $langs) {
    $langName = $lang.Name
    $fly      = $lang.Fly
    $ws = $wb.Worksheets.Item($langName)

    $targetDisplay   = "$uuid.md"
    $xlfDisplay      = "$uuid.$fileHash.$langName.xlf"

    $handoffUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($lang.HandoffSha)/ol-handoff/OpenLocalizationTestOrg/oltest-$fly-fly/xinjiang/ht/$xlfDisplay"
    $targetUrl   = "https://github.com/OpenLocalizationTestOrg/oltest-$fly-fly/blob/$($lang.TargetSha)/e2e/$targetDisplay"
    $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($lang.HandbackSha)/ol-handback/OpenLocalizationTestOrg/oltest-$fly-fly/xinjiang/ht/$xlfDisplay"

    # A: Source File Name (hyperlink to the .md source)
    $ws.Hyperlinks.Add($ws.Range("A4"), $overviewUrl, $null, $null, $targetDisplay)

    # B: File Extension
    $ws.Range("B4").Value = ".md"

    # C: Status
    $ws.Range("C4").Value = $status

    # D: Correspond Handoff File (hyperlink to handoff xlf)
    $ws.Hyperlinks.Add($ws.Range("D4"), $handoffUrl, $null, $null, $xlfDisplay)

    # E: Correspond Handoff Datetime
    $ws.Range("E4").Value = $lang.HandoffTime

    # F: Target File (hyperlink to the .md target)
    $ws.Hyperlinks.Add($ws.Range("F4"), $targetUrl, $null, $null, $targetDisplay)

    # G: Correspond Handback File (hyperlink to handback xlf)
    $ws.Hyperlinks.Add($ws.Range("G4"), $handbackUrl, $null, $null, $xlfDisplay)

    # H: Correspond Handback DateTime
    $ws.Range("H4").Value = $lang.HandbackTime

    # J: Handoff Reason
    $ws.Range("J4").Value = "Include"
}
